# Colors.xlsx - fill in the "Colors" (column H) description for the color
# combinations listed in rows 20-61. The shared-string table's new entries
# are created in the exact order the cells are written below (matching the
# uniqueCount/ordering seen in the target workbook), so the write order here
# intentionally does not follow row order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H43").Value = '백색;감색;청색;행황색;적색;두록색;흑청색'
$ws.Range("H47").Value = '구색;담자색;설백색;치색;자황색;백색;청색'
$ws.Range("H48").Value = '구색;담자색;백색;흑색;천청색;치색;회보라색;'
$ws.Range("H20").Value = '백색;초록색;구색;흑색;담황색;옥색;천청색'
$ws.Range("H21").Value = '황색;홍황색;담황색;연두색;자황색;백색;뇌록색'
$ws.Range("H22").Value = '담황색;연분홍색;백색;청벽색;회보라색;회색;비색'
$ws.Range("H31").Value = '백색;적색;청색;유황색;보라색;연두색;뇌록색'
$ws.Range("H32").Value = '황색;연두색;담황색;청색;백색;천청색;구색;흑록색'
$ws.Range("H33").Value = '담황색;옥색;진분홍색;청색;분홍색;적색'
$ws.Range("H34").Value = '옥색;비색;진초록색;연두색;흑록색;백색;청색'
$ws.Range("H35").Value = '담황색;자색;백색;흑색;분홍색;자황색;보라색;청색'
$ws.Range("H36").Value = '백색;적황색;소색;자황색;벽자색;청록색'
$ws.Range("H37").Value = '설백색;적색;백색;초록색;벽자색;소색;청색'
$ws.Range("H38").Value = '적색;청색;소색;구색;지황색;자색'
$ws.Range("H39").Value = '황색;적색;자주색;청색;춘유록색;연두색;담황색;보라색'
$ws.Range("H40").Value = '백색;자색;벽자색;흑색;회색;청색;적색'
$ws.Range("H41").Value = '적색;청색;적황색;담황색;자황색;황색;진분홍색'
$ws.Range("H42").Value = '백색;천청색;비색;흑색;구색;청색;연두색'
$ws.Range("H44").Value = '백색;흑색;자황색;양록색;적색;회보라색;벽청색'
$ws.Range("H45").Value = '구색;흑색;분홍색;청색;백색;자황색'
$ws.Range("H46").Value = '백색;흑색;토황색;석간주색;회색;홍람색;구색;청색'
$ws.Range("H59").Value = '담황색;연두색;분홍색;진분홍색;지황색;흑색;구색;청자색'
$ws.Range("H60").Value = '회보라색;남색;두록색;흑색;지황색;다자색;홍황색;양람색'
$ws.Range("H61").Value = '황색;벽청색'

# Leave the selection where the author's session ended up.
$ws.Range("H61").Select()
